# feat: add 2022-Q1 data
#
# The existing last sheet ("总计") becomes "2022-Q1" (keeps sheetId/rId),
# holding the same per-fund breakdown structure as the other quarterly
# sheets. A brand new "总计" sheet is appended after it with the updated
# roll-up table (2022-Q1 row added at the top, 2020-Q4 row re-added at the
# bottom).

$wb = $excel.ActiveWorkbook
$sheets = $wb.Worksheets

# ---------------------------------------------------------------------
# 1. Rename the current "总计" sheet -> "2022-Q1" and rebuild its content
# ---------------------------------------------------------------------
$q1 = $sheets.Item($sheets.Count)
$q1.Name = "2022-Q1"

# Extend the existing header style (already on B1:D1, s="2") across the
# new columns E1:H1 so every header cell shares one style.
$q1.Range("B1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Drop the old 5th data row (row 6) - the new sheet only has 4 data rows.
$q1.Rows.Item(6).Delete()

# Helper data: columns B..G must keep their original text representation
# (fund codes with leading zeros, decimal strings formatted to 2/4 dp),
# so every one of those is force-written as text. Column H is a genuine
# number (rank).
$q1Data = @(
    @(2, "539003", "建信富时100指数（QDII）人民币A", "0.71", "92.86", "5.71", "0.0405", 5),
    @(3, "008707", "建信富时100指数（QDII）美元现汇A", "0.71", "92.86", "5.71", "0.0405", 5),
    @(4, "008706", "建信富时100指数（QDII）人民币C", "0.20", "92.86", "5.71", "0.0114", 5),
    @(5, "008708", "建信富时100指数（QDII）美元现汇C", "0.20", "92.86", "5.71", "0.0114", 5)
)

foreach ($row in $q1Data) {
    $r = $row[0]
    $textCols = @("B", "C", "D", "E", "F", "G")
    for ($i = 0; $i -lt $textCols.Length; $i++) {
        $col = $textCols[$i]
        $cell = $q1.Range($col + $r)
        $cell.NumberFormat = "@"
        $cell.Value = $row[$i + 1]
        $cell.ClearFormats()
    }
    $q1.Range("H" + $r).Value = $row[7]
}

# ---------------------------------------------------------------------
# 2. Insert the new "总计" sheet right after "2022-Q1"
# ---------------------------------------------------------------------
$total = $sheets.Add([System.Reflection.Missing]::Value, $q1)
$total.Name = "总计"

# Reuse the same header/label style ("s=2") already resident in the
# workbook's style table instead of re-deriving one from scratch, so we
# don't spawn duplicate cellXfs entries.
$q1.Range("B1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$total.Range("A2:A7").PasteSpecial(-4122)

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$totalData = @(
    @(0, "2022-Q1", 4, 0.1),
    @(1, "2021-Q4", 4, 0.39),
    @(2, "2021-Q3", 4, 0.11),
    @(3, "2021-Q2", 4, 0.13),
    @(4, "2021-Q1", 4, 0.11),
    @(5, "2020-Q4", 4, 0.12)
)

foreach ($row in $totalData) {
    $r = $row[0] + 2
    $total.Range("A" + $r).Value = $row[0]
    $total.Range("B" + $r).Value = $row[1]
    $total.Range("C" + $r).Value = $row[2]
    $total.Range("D" + $r).Value = $row[3]
}
